$d = $word.ActiveDocument
$wXmlNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$wAllNs = $wXmlNs + ' xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'

# ---------------------------------------------------------------------------
# 1) Insert the two new "Direction's page" narrative paragraphs (plus the
#    blank spacer paragraph between them) right after the "aboutUs" bullet
#    and before the "reviews...?" bullet. The new paragraphs keep the
#    ListParagraph style but carry no bullet numbering.
# ---------------------------------------------------------------------------
$aboutUsPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "aboutUs*") {
        $aboutUsPara = $p
        break
    }
}

if ($aboutUsPara -ne $null) {
    # Create a clean, empty paragraph right after "aboutUs" to host the new XML.
    $aboutUsPara.Range.InsertParagraphAfter()

    $newHost = $aboutUsPara.Next()

    $newParasXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr><w:r><w:t>We&#8217;re proud to serve the men of Skippack, PA. Come on in and relax, knowing you&#8217;re in good hands. Whether you&#8217;re preparing for a big date night or prepping for your son&#8217;s Little League team photo, we&#8217;ll supply you with the cut you want.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr><w:r><w:t xml:space="preserve">We know that style is </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>personal</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> and a great cut accentuates one&#8217;s individuality. At Parker&#8217;s, we treat customers like friends and strive to deliver a truly customized experience for everyone. Come in, take a seat in our traditional, barber-styled chairs, lean back, and enjoy the experience of getting value for your dollar and cuts from those who enjoy their craft.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
'@

    $newHost.Range.InsertXML($newParasXml)
}

# ---------------------------------------------------------------------------
# 2) Drop the spell-check proofErr wrapper around "contactUs" (it was being
#    incorrectly flagged).
# ---------------------------------------------------------------------------
$contactUsPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "contactUs*") {
        $contactUsPara = $p
        break
    }
}

if ($contactUsPara -ne $null) {
    $cleanXml = '<w:p ' + $wAllNs + ' w14:paraId="3C2F365E" w14:textId="6AEFB8BD" w:rsidR="00F032E6" w:rsidRDefault="00F032E6" w:rsidP="00695C3F"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr></w:pPr><w:r><w:t>contactUs</w:t></w:r></w:p>'
    $contactUsPara.Range.InsertXML($cleanXml)
}

# ---------------------------------------------------------------------------
# 3) Remove the stray _GoBack bookmark that used to sit on the trailing empty
#    paragraph (it now lives on the new narrative paragraph instead).
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastCleanXml = '<w:p ' + $wAllNs + ' w14:paraId="3FA43555" w14:textId="77777777" w:rsidR="00CB5F3A" w:rsidRPr="00F45415" w:rsidRDefault="00CB5F3A" w:rsidP="00F12BDB"><w:pPr><w:pStyle w:val="ListParagraph"/><w:rPr><w:szCs w:val="24"/></w:rPr></w:pPr></w:p>'
$lastPara.Range.InsertXML($lastCleanXml)
